$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value2 = 159.25
$ws.Cells.Item(39, 9).Value2 = 76
$ws.Cells.Item(39, 10).Value2 = 218.71428
$ws.Cells.Item(39, 11).Value2 = 228
$ws.Cells.Item(39, 12).Value2 = 656.14284
$ws.Cells.Item(39, 13).Value2 = 68
$ws.Cells.Item(39, 14).Value2 = -1248.14284
$ws.Cells.Item(80, 8).Value2 = 1044.826
$ws.Cells.Item(80, 9).Value2 = 526
$ws.Cells.Item(80, 10).Value2 = 1443.9231
$ws.Cells.Item(80, 11).Value2 = 1578
$ws.Cells.Item(80, 12).Value2 = 4331.7693
$ws.Cells.Item(80, 13).Value2 = -580
$ws.Cells.Item(80, 14).Value2 = -6327.7693
$ws.Cells.Item(83, 8).Value2 = 1044.826
$ws.Cells.Item(83, 9).Value2 = 526
$ws.Cells.Item(83, 10).Value2 = 1443.9231
$ws.Cells.Item(83, 11).Value2 = 4734
$ws.Cells.Item(83, 12).Value2 = 12995.3079
$ws.Cells.Item(83, 13).Value2 = 258
$ws.Cells.Item(83, 14).Value2 = -22979.3079
$ws.Cells.Item(86, 8).Value2 = 2395.3076
$ws.Cells.Item(86, 9).Value2 = 2274.5
$ws.Cells.Item(86, 11).Value2 = 2274.5
$ws.Cells.Item(86, 13).Value2 = -1151.5
$ws.Cells.Item(89, 8).Value2 = 2395.3076
$ws.Cells.Item(89, 9).Value2 = 2274.5
$ws.Cells.Item(89, 11).Value2 = 11372.5
$ws.Cells.Item(89, 13).Value2 = -5756.5
$ws.Cells.Item(98, 8).Value2 = 2112.1538
$ws.Cells.Item(98, 9).Value2 = 2541.3
$ws.Cells.Item(98, 10).Value2 = 681.6667
$ws.Cells.Item(98, 11).Value2 = 2541.3
$ws.Cells.Item(98, 12).Value2 = 681.6667
$ws.Cells.Item(98, 13).Value2 = -1043.3
$ws.Cells.Item(98, 14).Value2 = -3677.6667
$ws.Cells.Item(122, 8).Value2 = 2112.1538
$ws.Cells.Item(122, 9).Value2 = 2541.3
$ws.Cells.Item(122, 10).Value2 = 681.6667
$ws.Cells.Item(122, 11).Value2 = 7623.900000000001
$ws.Cells.Item(122, 12).Value2 = 2045.0001
$ws.Cells.Item(122, 13).Value2 = -5173.900000000001
$ws.Cells.Item(122, 14).Value2 = -6945.0001
$ws.Cells.Item(132, 8).Value2 = 71433940
$ws.Cells.Item(132, 9).Value2 = 76928550
$ws.Cells.Item(132, 10).Value2 = 3999
$ws.Cells.Item(132, 11).Value2 = 230785650
$ws.Cells.Item(132, 12).Value2 = 11997
$ws.Cells.Item(132, 13).Value2 = -230783120
$ws.Cells.Item(132, 14).Value2 = -17057
$ws.Cells.Item(138, 8).Value2 = 2073.3774
$ws.Cells.Item(138, 9).Value2 = 1006.6829
$ws.Cells.Item(138, 10).Value2 = 5717.9165
$ws.Cells.Item(138, 11).Value2 = 3020.0487
$ws.Cells.Item(138, 12).Value2 = 17153.7495
$ws.Cells.Item(138, 13).Value2 = 2119.9513
$ws.Cells.Item(138, 14).Value2 = -27433.7495
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 4416.2095
$ws.Cells.Item(32, 9).Value2 = 3073.6858
$ws.Cells.Item(32, 11).Value2 = 3073.6858
$ws.Cells.Item(32, 13).Value2 = -2786.6858
$ws.Cells.Item(74, 8).Value2 = 50545.645
$ws.Cells.Item(74, 9).Value2 = 7037.6665
$ws.Cells.Item(74, 10).Value2 = 210074.89
$ws.Cells.Item(74, 11).Value2 = 7037.6665
$ws.Cells.Item(74, 12).Value2 = 210074.89
$ws.Cells.Item(74, 13).Value2 = -6163.6665
$ws.Cells.Item(74, 14).Value2 = -211822.89
$ws.Cells.Item(77, 8).Value2 = 50545.645
$ws.Cells.Item(77, 9).Value2 = 7037.6665
$ws.Cells.Item(77, 10).Value2 = 210074.89
$ws.Cells.Item(77, 11).Value2 = 35188.3325
$ws.Cells.Item(77, 12).Value2 = 1050374.45
$ws.Cells.Item(77, 13).Value2 = -30820.3325
$ws.Cells.Item(77, 14).Value2 = -1059110.45
$ws.Cells.Item(123, 8).Value2 = 55000
$ws.Cells.Item(123, 10).Value2 = 55000
$ws.Cells.Item(123, 12).Value2 = 55000
$ws.Cells.Item(123, 14).Value2 = -64800
$ws.Cells.Item(132, 8).Value2 = 2853.1538
$ws.Cells.Item(132, 9).Value2 = 2211.2354
$ws.Cells.Item(132, 10).Value2 = 3349.182
$ws.Cells.Item(132, 11).Value2 = 6633.706200000001
$ws.Cells.Item(132, 12).Value2 = 10047.546
$ws.Cells.Item(132, 13).Value2 = -4103.706200000001
$ws.Cells.Item(132, 14).Value2 = -15107.546
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value2 = 42000
$ws.Cells.Item(2, 10).Value2 = 42000
$ws.Cells.Item(2, 12).Value2 = 42000
$ws.Cells.Item(2, 14).Value2 = -42226
$ws.Cells.Item(99, 8).Value2 = 9527323
$ws.Cells.Item(99, 9).Value2 = 20410264
$ws.Cells.Item(99, 11).Value2 = 20410264
$ws.Cells.Item(99, 13).Value2 = -20408766
$ws.Cells.Item(107, 8).Value2 = 10209798
$ws.Cells.Item(107, 9).Value2 = 17860396
$ws.Cells.Item(107, 11).Value2 = 17860396
$ws.Cells.Item(107, 13).Value2 = -17858476
$ws.Cells.Item(134, 8).Value2 = 4631.6
$ws.Cells.Item(134, 9).Value2 = 2353.2
$ws.Cells.Item(134, 10).Value2 = 8049.2
$ws.Cells.Item(134, 11).Value2 = 7059.599999999999
$ws.Cells.Item(134, 12).Value2 = 24147.6
$ws.Cells.Item(134, 13).Value2 = -4524.599999999999
$ws.Cells.Item(134, 14).Value2 = -29217.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value2 = 42598.816
$ws.Cells.Item(132, 9).Value2 = 2299.3333
$ws.Cells.Item(132, 11).Value2 = 6897.999899999999
$ws.Cells.Item(132, 13).Value2 = -4367.999899999999
$ws.Cells.Item(134, 8).Value2 = 3082.476
$ws.Cells.Item(134, 9).Value2 = 2229.4167
$ws.Cells.Item(134, 11).Value2 = 6688.250100000001
$ws.Cells.Item(134, 13).Value2 = -4153.250100000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value2 = 514.1111
$ws.Cells.Item(47, 9).Value2 = 515.875
$ws.Cells.Item(47, 11).Value2 = 1547.625
$ws.Cells.Item(47, 13).Value2 = -1116.625
$ws.Cells.Item(69, 8).Value2 = 4292.3335
$ws.Cells.Item(69, 9).Value2 = 1011
$ws.Cells.Item(69, 10).Value2 = 4948.6
$ws.Cells.Item(69, 11).Value2 = 3033
$ws.Cells.Item(69, 12).Value2 = 14845.8
$ws.Cells.Item(69, 13).Value2 = -2222
$ws.Cells.Item(69, 14).Value2 = -16467.8
$ws.Cells.Item(72, 8).Value2 = 4292.3335
$ws.Cells.Item(72, 9).Value2 = 1011
$ws.Cells.Item(72, 10).Value2 = 4948.6
$ws.Cells.Item(72, 11).Value2 = 9099
$ws.Cells.Item(72, 12).Value2 = 44537.4
$ws.Cells.Item(72, 13).Value2 = -5043
$ws.Cells.Item(72, 14).Value2 = -52649.4
$ws.Cells.Item(87, 8).Value2 = 9399.6
$ws.Cells.Item(87, 9).Value2 = 9399.6
$ws.Cells.Item(87, 11).Value2 = 28198.8
$ws.Cells.Item(87, 13).Value2 = -26950.8
$ws.Cells.Item(90, 8).Value2 = 9399.6
$ws.Cells.Item(90, 9).Value2 = 9399.6
$ws.Cells.Item(90, 11).Value2 = 84596.40000000001
$ws.Cells.Item(90, 13).Value2 = -78356.40000000001
$ws.Cells.Item(107, 8).Value2 = 1669.6428
$ws.Cells.Item(107, 10).Value2 = 2099.7
$ws.Cells.Item(107, 12).Value2 = 6299.099999999999
$ws.Cells.Item(107, 14).Value2 = -10139.1
$ws.Cells.Item(114, 8).Value2 = 468
$ws.Cells.Item(114, 9).Value2 = 485.25
$ws.Cells.Item(114, 10).Value2 = 445
$ws.Cells.Item(114, 11).Value2 = 1455.75
$ws.Cells.Item(114, 12).Value2 = 1335
$ws.Cells.Item(114, 13).Value2 = 1798.25
$ws.Cells.Item(114, 14).Value2 = -7843
$ws.Cells.Item(124, 8).Value2 = 2400
$ws.Cells.Item(124, 9).Value2 = 2400
$ws.Cells.Item(124, 11).Value2 = 7200
$ws.Cells.Item(124, 13).Value2 = -2290
$ws.Cells.Item(128, 8).Value2 = 181750
$ws.Cells.Item(128, 9).Value2 = 181750
$ws.Cells.Item(128, 11).Value2 = 545250
$ws.Cells.Item(128, 13).Value2 = -540270
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value2 = 1725.5
$ws.Cells.Item(107, 9).Value2 = 1700.6666
$ws.Cells.Item(107, 10).Value2 = 1800
$ws.Cells.Item(107, 11).Value2 = 1700.6666
$ws.Cells.Item(107, 12).Value2 = 1800
$ws.Cells.Item(107, 13).Value2 = 219.3334
$ws.Cells.Item(107, 14).Value2 = -5640
$ws.Cells.Item(122, 8).Value2 = 359785.44
$ws.Cells.Item(122, 10).Value2 = 6105.2856
$ws.Cells.Item(122, 12).Value2 = 18315.8568
$ws.Cells.Item(122, 14).Value2 = -23215.8568
$ws.Cells.Item(132, 8).Value2 = 3552.75
$ws.Cells.Item(132, 9).Value2 = 2946.4167
$ws.Cells.Item(132, 10).Value2 = 4765.4165
$ws.Cells.Item(132, 11).Value2 = 8839.250100000001
$ws.Cells.Item(132, 12).Value2 = 14296.2495
$ws.Cells.Item(132, 13).Value2 = -6309.250100000001
$ws.Cells.Item(132, 14).Value2 = -19356.2495
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 3551.95
$ws.Cells.Item(7, 9).Value2 = 1389.2
$ws.Cells.Item(7, 11).Value2 = 1389.2
$ws.Cells.Item(7, 13).Value2 = -1277.2
$ws.Cells.Item(22, 8).Value2 = 60625.465
$ws.Cells.Item(22, 9).Value2 = 296829.34
$ws.Cells.Item(22, 10).Value2 = 1574.5
$ws.Cells.Item(22, 11).Value2 = 296829.34
$ws.Cells.Item(22, 12).Value2 = 1574.5
$ws.Cells.Item(22, 13).Value2 = -296534.34
$ws.Cells.Item(22, 14).Value2 = -2164.5
$ws.Cells.Item(27, 8).Value2 = 60625.465
$ws.Cells.Item(27, 9).Value2 = 296829.34
$ws.Cells.Item(27, 10).Value2 = 1574.5
$ws.Cells.Item(27, 11).Value2 = 296829.34
$ws.Cells.Item(27, 12).Value2 = 1574.5
$ws.Cells.Item(27, 13).Value2 = -296722.34
$ws.Cells.Item(27, 14).Value2 = -1788.5
$ws.Cells.Item(61, 8).Value2 = 3970510
$ws.Cells.Item(61, 9).Value2 = 4632118.5
$ws.Cells.Item(61, 10).Value2 = 859.75
$ws.Cells.Item(61, 11).Value2 = 4632118.5
$ws.Cells.Item(61, 12).Value2 = 859.75
$ws.Cells.Item(61, 13).Value2 = -4631916.5
$ws.Cells.Item(61, 14).Value2 = -1263.75
$ws.Cells.Item(107, 8).Value2 = 3089
$ws.Cells.Item(107, 9).Value2 = 3089
$ws.Cells.Item(107, 11).Value2 = 3089
$ws.Cells.Item(107, 13).Value2 = -1169
$ws.Cells.Item(113, 8).Value2 = 3970510
$ws.Cells.Item(113, 9).Value2 = 4632118.5
$ws.Cells.Item(113, 10).Value2 = 859.75
$ws.Cells.Item(113, 11).Value2 = 4632118.5
$ws.Cells.Item(113, 12).Value2 = 859.75
$ws.Cells.Item(113, 13).Value2 = -4629948.5
$ws.Cells.Item(113, 14).Value2 = -5199.75
$ws.Cells.Item(126, 8).Value2 = 3551.95
$ws.Cells.Item(126, 9).Value2 = 1389.2
$ws.Cells.Item(126, 11).Value2 = 4167.6
$ws.Cells.Item(126, 13).Value2 = -1697.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(106, 8).Value2 = 44521
$ws.Cells.Item(106, 9).Value2 = 20342
$ws.Cells.Item(106, 10).Value2 = 68700
$ws.Cells.Item(106, 11).Value2 = 20342
$ws.Cells.Item(106, 12).Value2 = 68700
$ws.Cells.Item(106, 13).Value2 = -19080
$ws.Cells.Item(106, 14).Value2 = -71224
$ws.Cells.Item(114, 8).Value2 = 49982
$ws.Cells.Item(114, 10).Value2 = 49982
$ws.Cells.Item(114, 12).Value2 = 49982
$ws.Cells.Item(114, 14).Value2 = -58660
$ws.Cells.Item(132, 8).Value2 = 32292296
$ws.Cells.Item(132, 9).Value2 = 52639176
$ws.Cells.Item(132, 10).Value2 = 76401.336
$ws.Cells.Item(132, 11).Value2 = 157917528
$ws.Cells.Item(132, 12).Value2 = 229204.008
$ws.Cells.Item(132, 13).Value2 = -157914998
$ws.Cells.Item(132, 14).Value2 = -234264.008
